$wb = $excel.ActiveWorkbook

# "NextTransaction.xaml deleted and replaced with actions in transitions."
# Remove the now-obsolete "Next" workblock settings (wbNextTransaction_Type /
# wbNextTransaction_SuppressSuccessful) from the Workblocks sheet - rows 11:12.
$wsWorkblocks = $wb.Worksheets.Item("Workblocks")
$wsWorkblocks.Rows("11:12").Delete()

# Update the remembered selection on the Workblocks sheet.
$wsWorkblocks.Activate()
$wsWorkblocks.Range("C24").Select()

# Make the Tasks sheet the active / selected tab (was Settings before).
$wsTasks = $wb.Worksheets.Item("Tasks")
$wsTasks.Activate()

Write-Host "Removed NextTransaction workblock settings and updated active sheet/selection"
